$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("HICOOL ELECTRONIC INDUSTRIES", "MUMBAI", "MAHARASHTRA"),
    @("HCE DYNAMICS PRIVATE LIMITED", "MUMBAI", "MAHARASHTRA"),
    @("PARKASH ELECTRIC COMPANY", "CHANDIGARH", "CHANDIGARH"),
    @("N-RACK ACCESSORIES PVT LTD", "BANGALORE", "KARNATAKA")
)

$startRow = 274
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
